$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.379.99'
$ws.Range("E2").Value = '  +0.79%  '
$ws.Range("D3").Value = '3.802.39'
$ws.Range("E3").Value = '  -0.36%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '''699.04'
$ws.Range("E5").Value = '  +5.38%  '
$ws.Range("D6").Value = '''174.05'
$ws.Range("E6").Value = '  +2.90%  '
$ws.Range("D7").Value = '3.802.90'
$ws.Range("E7").Value = '  -0.24%  '
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").Value = '''0.528'
$ws.Range("E9").Value = '  +0.01%  '
$ws.Range("E10").Value = '  +0.78%  '
$ws.Range("D11").Value = '''7.42'
$ws.Range("E11").Value = '  +6.33%  '
$ws.Range("E12").Value = '  -0.06%  '
$ws.Range("D13").Value = '''0.0000257'
$ws.Range("E13").Value = '  +5.17%  '
$ws.Range("D14").Value = '''36.59'
$ws.Range("E14").Value = '  +2.31%  '
$ws.Range("D15").Value = '4.446.55'
$ws.Range("E15").Value = '  -0.34%  '
$ws.Range("D16").Value = '3.805.90'
$ws.Range("E16").Value = '  +0.18%  '
$ws.Range("D17").Value = '71.358.67'
$ws.Range("E17").Value = '  +0.76%  '
$ws.Range("E18").Value = '  -0.47%  '
$ws.Range("E19").Value = '  +0.88%  '
$ws.Range("E20").Value = '  +0.35%  '
$ws.Range("D21").Value = '''11.13'
$ws.Range("E21").Value = '  +7.98%  '
$ws.Range("D22").Value = '''484.75'
$ws.Range("E22").Value = '  +1.21%  '
$ws.Range("E23").Value = '  +0.38%  '
$ws.Range("D24").Value = '''84.49'
$ws.Range("E24").Value = '  +1.84%  '
$ws.Range("D25").Value = '''0.0000142'
$ws.Range("E25").Value = '  -2.93%  '
$ws.Range("E26").Value = '  +0.32%  '
$ws.Range("D27").Value = '''10.53'
$ws.Range("E27").Value = '  +1.40%  '
$ws.Range("E28").Value = '  +1.38%  '
$ws.Range("D29").Value = '3.953.86'
$ws.Range("E29").Value = '  -0.37%  '
$ws.Range("E30").Value = '  -0.04%  '
$ws.Range("D31").Value = '''3.12'
$ws.Range("E31").Value = '  +9.98%  '
$ws.Range("D32").Value = '''2.31'
$ws.Range("E32").Value = '  -0.25%  '
$ws.Range("E33").Value = '  +1.59%  '
$ws.Range("D34").Value = '''0.184'
$ws.Range("E34").Value = '  +2.59%  '
$ws.Range("D35").Value = '''29.58'
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("D36").Value = '''9.29'
$ws.Range("E36").Value = '  +2.20%  '
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("D38").Value = '''0.103'
$ws.Range("E38").Value = '  +1.45%  '
$ws.Range("E39").Value = '  +13.36%  '
$ws.Range("D40").Value = '''3.41'
$ws.Range("E40").Value = '  -0.67%  '
$ws.Range("D41").Value = '''6.03'
$ws.Range("E41").Value = '  +1.84%  '
$ws.Range("D42").Value = '''0.996'
$ws.Range("E42").Value = '  +2.77%  '
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("D45").Value = '''164.71'
$ws.Range("E45").Value = '  +3.72%  '
$ws.Range("D46").Value = '''0.000305'
$ws.Range("E46").Value = '  +5.02%  '
$ws.Range("D47").Value = '''44.78'
$ws.Range("E47").Value = '  -1.31%  '
$ws.Range("D48").Value = '''48.58'
$ws.Range("E48").Value = '  +0.37%  '
$ws.Range("E49").Value = '  +0.73%  '
$ws.Range("D50").Value = '''416.63'
$ws.Range("E50").Value = '  +4.93%  '
$ws.Range("D51").Value = '''8.66'
$ws.Range("E51").Value = '  +1.92%  '

# The cells above hold numeric-looking text (e.g. "699.04") in the source
# data; a leading apostrophe keeps Excel from silently coercing them to
# floating point numbers. Re-apply the default "Normal" cell style so the
# apostrophe's implicit quote-prefix formatting doesn't linger on the cell.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"

Write-Output "Updated cryptos list"
